# Scheduled-runner refresh of Zalera_Profits pricing data.
# Updates currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns
# (H:N) for the rows whose underlying market prices moved, across the
# ALC / ARM / BSM / CRP / CUL / LTW / WVR sheets (GSM unchanged).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 961632.25
$ws.Range("I9").Value = 1576.7
$ws.Range("J9").Value = 1601669.4
$ws.Range("K9").Value = 1576.7
$ws.Range("L9").Value = 1601669.4
$ws.Range("M9").Value = -1407.7
$ws.Range("N9").Value = -1602007.4
$ws.Range("H40").Value = 1829.9535
$ws.Range("I40").Value = 1787.25
$ws.Range("K40").Value = 1787.25
$ws.Range("M40").Value = -1612.25
$ws.Range("H88").Value = 10499.714
$ws.Range("J88").Value = 9166.333000000001
$ws.Range("L88").Value = 9166.333000000001
$ws.Range("N88").Value = -9978.333000000001
$ws.Range("H91").Value = 10499.714
$ws.Range("J91").Value = 9166.333000000001
$ws.Range("L91").Value = 9166.333000000001
$ws.Range("N91").Value = -11974.333
$ws.Range("H116").Value = 6499.9
$ws.Range("I116").Value = 5625
$ws.Range("J116").Value = 9999.5
$ws.Range("K116").Value = 5625
$ws.Range("L116").Value = 9999.5
$ws.Range("M116").Value = -2183
$ws.Range("N116").Value = -16883.5
$ws.Range("H129").Value = 1722.1904
$ws.Range("I129").Value = 960.5454999999999
$ws.Range("J129").Value = 2560
$ws.Range("K129").Value = 2881.6365
$ws.Range("L129").Value = 7680
$ws.Range("M129").Value = 2118.3635
$ws.Range("N129").Value = -17680
$ws.Range("H132").Value = 1709.8286
$ws.Range("I132").Value = 1093.421
$ws.Range("J132").Value = 2441.8125
$ws.Range("K132").Value = 3280.263
$ws.Range("L132").Value = 7325.4375
$ws.Range("M132").Value = -750.2629999999999
$ws.Range("N132").Value = -12385.4375
$ws.Range("H137").Value = 9147.333000000001
$ws.Range("I137").Value = 6312.5
$ws.Range("J137").Value = 10178.182
$ws.Range("K137").Value = 18937.5
$ws.Range("L137").Value = 30534.546
$ws.Range("M137").Value = -16387.5
$ws.Range("N137").Value = -35634.546
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3468134.5
$ws.Range("I32").Value = 3656572.8
$ws.Range("K32").Value = 3656572.8
$ws.Range("M32").Value = -3656285.8
$ws.Range("H61").Value = 4284.0415
$ws.Range("I61").Value = 3686.524
$ws.Range("J61").Value = 8466.666999999999
$ws.Range("K61").Value = 3686.524
$ws.Range("L61").Value = 8466.666999999999
$ws.Range("M61").Value = -3474.524
$ws.Range("N61").Value = -8890.666999999999
$ws.Range("H74").Value = 378800.06
$ws.Range("I74").Value = 533173.9
$ws.Range("K74").Value = 533173.9
$ws.Range("M74").Value = -532299.9
$ws.Range("H77").Value = 378800.06
$ws.Range("I77").Value = 533173.9
$ws.Range("K77").Value = 2665869.5
$ws.Range("M77").Value = -2661501.5
$ws.Range("H132").Value = 8731
$ws.Range("I132").Value = 6187.933
$ws.Range("J132").Value = 13499.25
$ws.Range("K132").Value = 18563.799
$ws.Range("L132").Value = 40497.75
$ws.Range("M132").Value = -16033.799
$ws.Range("N132").Value = -45557.75
$ws.Range("H136").Value = 4284.0415
$ws.Range("I136").Value = 3686.524
$ws.Range("J136").Value = 8466.666999999999
$ws.Range("K136").Value = 11059.572
$ws.Range("L136").Value = 25400.001
$ws.Range("M136").Value = -8509.572
$ws.Range("N136").Value = -30500.001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 428000
$ws.Range("I96").Value = 428000
$ws.Range("K96").Value = 428000
$ws.Range("M96").Value = -425254
$ws.Range("H134").Value = 6251.3
$ws.Range("I134").Value = 4127.25
$ws.Range("K134").Value = 12381.75
$ws.Range("M134").Value = -9846.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4174.421
$ws.Range("I31").Value = 2133.3076
$ws.Range("J31").Value = 5235.8
$ws.Range("K31").Value = 2133.3076
$ws.Range("L31").Value = 5235.8
$ws.Range("M31").Value = -1838.3076
$ws.Range("N31").Value = -5825.8
$ws.Range("H34").Value = 4174.421
$ws.Range("I34").Value = 2133.3076
$ws.Range("J34").Value = 5235.8
$ws.Range("K34").Value = 2133.3076
$ws.Range("L34").Value = 5235.8
$ws.Range("M34").Value = -1931.3076
$ws.Range("N34").Value = -5639.8
$ws.Range("H58").Value = 8248.25
$ws.Range("I58").Value = 7104.643
$ws.Range("J58").Value = 8864.038
$ws.Range("K58").Value = 7104.643
$ws.Range("L58").Value = 8864.038
$ws.Range("M58").Value = -6901.643
$ws.Range("N58").Value = -9270.038
$ws.Range("H105").Value = 1738.5
$ws.Range("I105").Value = 1738.5
$ws.Range("K105").Value = 1738.5
$ws.Range("M105").Value = 8.5
$ws.Range("H122").Value = 168720.5
$ws.Range("I122").Value = 251251.75
$ws.Range("J122").Value = 3658
$ws.Range("K122").Value = 753755.25
$ws.Range("L122").Value = 10974
$ws.Range("M122").Value = -751305.25
$ws.Range("N122").Value = -15874
$ws.Range("H134").Value = 8731.941000000001
$ws.Range("I134").Value = 8792.933999999999
$ws.Range("J134").Value = 8274.5
$ws.Range("K134").Value = 26378.802
$ws.Range("L134").Value = 24823.5
$ws.Range("M134").Value = -23843.802
$ws.Range("N134").Value = -29893.5
$ws.Range("H136").Value = 8248.25
$ws.Range("I136").Value = 7104.643
$ws.Range("J136").Value = 8864.038
$ws.Range("K136").Value = 21313.929
$ws.Range("L136").Value = 26592.114
$ws.Range("M136").Value = -18763.929
$ws.Range("N136").Value = -31692.114
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 42.714287
$ws.Range("I23").Value = 29.5
$ws.Range("J23").Value = 48
$ws.Range("K23").Value = 88.5
$ws.Range("L23").Value = 144
$ws.Range("M23").Value = 146.5
$ws.Range("N23").Value = -614
$ws.Range("H68").Value = 93941.74000000001
$ws.Range("I68").Value = 288445.44
$ws.Range("J68").Value = 8846.375
$ws.Range("K68").Value = 865336.3200000001
$ws.Range("L68").Value = 26539.125
$ws.Range("M68").Value = -864525.3200000001
$ws.Range("N68").Value = -28161.125
$ws.Range("H71").Value = 93941.74000000001
$ws.Range("I71").Value = 288445.44
$ws.Range("J71").Value = 8846.375
$ws.Range("K71").Value = 2596008.96
$ws.Range("L71").Value = 79617.375
$ws.Range("M71").Value = -2591952.96
$ws.Range("N71").Value = -87729.375
$ws.Range("H107").Value = 5213.846
$ws.Range("I107").Value = 611.4286
$ws.Range("J107").Value = 10583.333
$ws.Range("K107").Value = 1834.2858
$ws.Range("L107").Value = 31749.999
$ws.Range("M107").Value = 85.71420000000012
$ws.Range("N107").Value = -35589.999
$ws.Range("H129").Value = 26316634
$ws.Range("I129").Value = 791.7273
$ws.Range("J129").Value = 62500920
$ws.Range("K129").Value = 2375.1819
$ws.Range("L129").Value = 187502760
$ws.Range("M129").Value = 2624.8181
$ws.Range("N129").Value = -187512760
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4769004.5
$ws.Range("I46").Value = 25002250
$ws.Range("J46").Value = 8240.647000000001
$ws.Range("K46").Value = 25002250
$ws.Range("L46").Value = 8240.647000000001
$ws.Range("M46").Value = -25002062
$ws.Range("N46").Value = -8616.647000000001
$ws.Range("H55").Value = 797.6923
$ws.Range("I55").Value = 247.125
$ws.Range("J55").Value = 1678.6
$ws.Range("K55").Value = 247.125
$ws.Range("L55").Value = 1678.6
$ws.Range("M55").Value = -74.125
$ws.Range("N55").Value = -2024.6
$ws.Range("H132").Value = 5607.515
$ws.Range("I132").Value = 3401.76
$ws.Range("J132").Value = 12500.5
$ws.Range("K132").Value = 10205.28
$ws.Range("L132").Value = 37501.5
$ws.Range("M132").Value = -7675.280000000001
$ws.Range("N132").Value = -42561.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2713.32
$ws.Range("I100").Value = 3286.9473
$ws.Range("J100").Value = 896.8333
$ws.Range("K100").Value = 6573.8946
$ws.Range("L100").Value = 1793.6666
$ws.Range("M100").Value = -6032.8946
$ws.Range("N100").Value = -2875.6666
$ws.Range("H126").Value = 25252224
$ws.Range("I126").Value = 25252224
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 75756672
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -75754202
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 7112.0605
$ws.Range("I132").Value = 5588.1763
$ws.Range("J132").Value = 8731.1875
$ws.Range("K132").Value = 16764.5289
$ws.Range("L132").Value = 26193.5625
$ws.Range("M132").Value = -14234.5289
$ws.Range("N132").Value = -31253.5625
